# Mon 29 Aug 2016 21:38 commit by king
# Functionality : Added GWT Server & User checking support.
# Issue Fix : Edit Regular/Special Schedule not working without changing the name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new "User checking" rows under the
#     "GWT to Server Packet Format (SSL)" table (was rows 18-19, becomes 18-20)
#     and under the "Server to GWT Packet Format (SSL)" table (was row 22, becomes 24-26).
$ws.Rows("19:20").Insert()
$ws.Rows("25:26").Insert()

# --- GWT to Server Packet Format (SSL) table (header row stays at 17) ---
# Row 18: UserCheckNameExists(1, Username)  -- replaces old SQL Statement row
$ws.Range("A18").Value = "UserCheckNameExists"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "Username"
$ws.Range("G18").ClearContents()
$ws.Range("H18").ClearContents()

# Row 19 (new): UserRegister(2, Username, SHA-1(Password))
$ws.Range("A19").Value = "UserRegister"
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Username"
$ws.Range("D19").Value = "SHA-1(Password)"

# Row 20 (new): UserCheckCredentialOK(3, Username, SHA-1(Password))
$ws.Range("A20").Value = "UserCheckCredentialOK"
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "Username"
$ws.Range("D20").Value = "SHA-1(Password)"

# --- Server to GWT Packet Format (SSL) table (header text at A22, column headers at row 23) ---
# Row 24: UserCheckNameExists -> boolean -- replaces old SQL Query Result row
$ws.Range("A24").Value = "UserCheckNameExists"
$ws.Range("B24").ClearContents()
$ws.Range("C24").Value = "boolean"
$ws.Range("G24").Value = "To be done"
$ws.Range("H24").Value = "v1.0"

# Row 25 (new): UserRegister -> boolean
$ws.Range("A25").Value = "UserRegister"
$ws.Range("C25").Value = "boolean"

# Row 26 (new): UserCheckCredentialOK -> FAILCHECK/ACTIVE/PENDING APPROVAL/DEACTIVATED
$ws.Range("A26").Value = "UserCheckCredentialOK"
$ws.Range("C26").Value = "FAILCHECK/ACTIVE/PENDING APPROVAL/DEACTIVATED"

# --- Update the view: scrolled down a bit, with C30 as the active cell ---
$ws.Range("C30").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
